# Weekly CompStat refresh: new crime data collected.
# Applies the week-over-week numeric refresh plus the header volume/date bump.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: force a cell that currently holds a NUMBER to instead hold the
# literal text "0" (matching the report's "no activity" placeholder cells),
# while keeping its General number format (mirrors the other "0"-placeholder
# cells such as C14/D14 elsewhere in the sheet).
# ---------------------------------------------------------------------------
function Set-ZeroPlaceholderText($addr) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = "0"
    $ws.Range($addr).NumberFormat = "General"
}

# ---------------------------------------------------------------------------
# Helper: force a cell that currently holds the literal text "0" to instead
# hold a real number, restoring the normal integer number format used by
# sibling count cells.
# ---------------------------------------------------------------------------
function Set-NumberFromPlaceholder($addr, $value) {
    $ws.Range($addr).Value = $value
    $ws.Range($addr).NumberFormat = "#,##0"
}

# ---------------------------------------------------------------------------
# Header: volume number 35 -> 36, and week-covering dates roll forward a week.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/1/2025  Through  9/7/2025"

# ---------------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------------
Set-NumberFromPlaceholder "I14" 1
$ws.Range("K14").Value = -85.714285714285
$ws.Range("L14").Value = -83.333333333333
$ws.Range("M14").Value = -87.5
$ws.Range("N14").Value = -96.153846153846

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
Set-ZeroPlaceholderText "C15"
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = 40
$ws.Range("M15").Value = 7.692307692307
$ws.Range("N15").Value = -50

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
Set-NumberFromPlaceholder "C16" 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = -57.894736842105
$ws.Range("I16").Value = 131
$ws.Range("J16").Value = 178
$ws.Range("K16").Value = -26.404494382022
$ws.Range("L16").Value = -28.021978021978
$ws.Range("M16").Value = -13.815789473684
$ws.Range("N16").Value = -77.871621621621

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 37
$ws.Range("H17").Value = -27.027027027027
$ws.Range("I17").Value = 240
$ws.Range("J17").Value = 344
$ws.Range("K17").Value = -30.232558139534
$ws.Range("L17").Value = -11.764705882352
$ws.Range("M17").Value = 48.148148148148
$ws.Range("N17").Value = -48.164146868250

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 33.333333333333
$ws.Range("I18").Value = 99
$ws.Range("J18").Value = 69
$ws.Range("K18").Value = 43.478260869565
$ws.Range("L18").Value = 22.222222222222
$ws.Range("M18").Value = 28.571428571428
$ws.Range("N18").Value = -71.633237822349

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 5
$ws.Range("G19").Value = 40
$ws.Range("H19").Value = -32.5
$ws.Range("I19").Value = 275
$ws.Range("J19").Value = 358
$ws.Range("K19").Value = -23.184357541899
$ws.Range("L19").Value = -31.421446384039
$ws.Range("M19").Value = 78.571428571428
$ws.Range("N19").Value = -1.433691756272

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
Set-NumberFromPlaceholder "C20" 1
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 50
$ws.Range("J20").Value = 79
$ws.Range("K20").Value = -36.708860759493
$ws.Range("L20").Value = -15.254237288135
$ws.Range("M20").Value = 2.040816326530
$ws.Range("N20").Value = -78.354978354978

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 18
$ws.Range("E21").Value = -25
$ws.Range("F21").Value = 77
$ws.Range("G21").Value = 108
$ws.Range("H21").Value = -28.703703703703
$ws.Range("I21").Value = 810
$ws.Range("J21").Value = 1045
$ws.Range("K21").Value = -22.488038277512
$ws.Range("L21").Value = -19.722497522299
$ws.Range("M21").Value = 31.707317073170
$ws.Range("N21").Value = -58.841463414634

# ---------------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------------
$ws.Range("F22").Value = 4
$ws.Range("I22").Value = 17
$ws.Range("K22").Value = 30.769230769230
$ws.Range("L22").Value = 6.25
$ws.Range("M22").Value = -39.285714285714

# ---------------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------------
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 12
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = 9.090909090909
$ws.Range("I23").Value = 121
$ws.Range("J23").Value = 116
$ws.Range("K23").Value = 4.310344827586
$ws.Range("L23").Value = 4.310344827586
$ws.Range("M23").Value = 68.055555555555

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 73.333333333333
$ws.Range("F24").Value = 75
$ws.Range("G24").Value = 67
$ws.Range("H24").Value = 11.940298507462
$ws.Range("I24").Value = 556
$ws.Range("J24").Value = 691
$ws.Range("K24").Value = -19.536903039073
$ws.Range("L24").Value = -35.498839907192
$ws.Range("M24").Value = -13.125

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 300
$ws.Range("F25").Value = 16
$ws.Range("G25").Value = 8
$ws.Range("H25").Value = 100
$ws.Range("I25").Value = 130
$ws.Range("J25").Value = 136
$ws.Range("K25").Value = -4.411764705882
$ws.Range("L25").Value = -69.047619047619

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = -28.571428571428
$ws.Range("F26").Value = 35
$ws.Range("H26").Value = -50
$ws.Range("I26").Value = 434
$ws.Range("J26").Value = 538
$ws.Range("K26").Value = -19.330855018587
$ws.Range("L26").Value = 23.646723646723
$ws.Range("M26").Value = 26.162790697674

# ---------------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------------
Set-ZeroPlaceholderText "C27"
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -75
$ws.Range("J27").Value = 15
$ws.Range("K27").Value = 0

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 40
$ws.Range("I28").Value = 50
$ws.Range("J28").Value = 43
$ws.Range("K28").Value = 16.279069767441
$ws.Range("L28").Value = 28.205128205128

# ---------------------------------------------------------------------------
# Row 29 - Shooting Vic.
# ---------------------------------------------------------------------------
$ws.Range("M29").Value = -75.862068965517
$ws.Range("N29").Value = -85.416666666666

# ---------------------------------------------------------------------------
# Row 30 - Shooting Inc.
# ---------------------------------------------------------------------------
$ws.Range("M30").Value = -79.166666666666
$ws.Range("N30").Value = -89.361702127659
